$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen new column F to fit the source URLs
$ws.Columns("F").ColumnWidth = 71.16

# New "source" header in F1, matching the style of the other header cells (A1:E1)
$ws.Range("F1").Value2 = "source"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Freesound.org source URLs for rows 3-15, inserted as hyperlinks (display text = URL)
$urls = @(
  "https://freesound.org/people/Disagree/sounds/433725/",
  "https://freesound.org/people/spycrah/sounds/471097/",
  "https://freesound.org/people/SoundFlakes/sounds/492239/",
  "https://freesound.org/people/qubodup/sounds/60027/",
  "https://freesound.org/people/dereklieu/sounds/241822/",
  "https://freesound.org/people/xpoki/sounds/432755/",
  "https://freesound.org/people/Michel88/sounds/76959/",
  "https://freesound.org/people/Aleks41/sounds/449552/",
  "https://freesound.org/people/238310/sounds/370189/",
  "https://freesound.org/people/Meisben/sounds/488037/",
  "https://freesound.org/people/JoelAudio/sounds/135463/",
  "https://freesound.org/people/cmorris035/sounds/319152/",
  "https://freesound.org/people/FreqMan/sounds/23168/"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
  $row = $i + 3
  $cell = $ws.Range("F$row")
  $ws.Hyperlinks.Add($cell, $urls[$i]) | Out-Null
}

$ws.Range("F16").Select() | Out-Null
